# Macroferia Regional de Talca - Haba
# A new weekly price record is inserted as row 48, pushing the existing
# rows 48-109 down to 49-110 (the sheet simply gets one more data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; this shifts rows 48..109 down to 49..110
# and keeps the inherited formatting (e.g. the date style on column D).
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly record.
$ws.Cells.Item(48, 1).Value  = 5
$ws.Cells.Item(48, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(48, 3).Value  = "Maule"
$ws.Cells.Item(48, 4).Value  = 44880
$ws.Cells.Item(48, 5).Value  = 7
$ws.Cells.Item(48, 6).Value  = 100112026
$ws.Cells.Item(48, 7).Value  = "Haba"
$ws.Cells.Item(48, 8).Value  = "Sin especificar"
$ws.Cells.Item(48, 9).Value  = "Primera"
$ws.Cells.Item(48, 10).Value = 400
$ws.Cells.Item(48, 11).Value = 9000
$ws.Cells.Item(48, 12).Value = 9000
$ws.Cells.Item(48, 13).Value = 9000
$ws.Cells.Item(48, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(48, 15).Value = "Región del Maule"
$ws.Cells.Item(48, 16).Value = 360
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
